# Exclude denied embassy requests
# Removes the row for "The beforelife" (row 68), which shifts the rows
# below it up by one (Ben -> row 68, Brest Oblast -> row 69), and shrinks
# the used range from A1:I70 to A1:I69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(68).Delete()
